$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

# Update source citation block (B3:B7)
$about.Range("B3").Value = "Massachusetts Institute of Technology"
$about.Range("B4").Value = 2021
$about.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$about.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$about.Range("B7").Value = "Abstract"

# Clear the old footnote that referenced the retired chart image
$about.Range("C8").Value = ""

# Add the new footnote about the averaged learning rate
$about.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Remove the now-unused chart picture from the About sheet
if ($about.Shapes.Count -gt 0) {
    for ($i = $about.Shapes.Count; $i -ge 1; $i--) {
        $about.Shapes.Item($i).Delete()
    }
}

# --- "PDiBCpDoC" sheet updates ---
$data = $wb.Worksheets.Item("PDiBCpDoC")

# Replace the hard-coded decline rate with the average of the two cited rates
$data.Range("B2").Formula = "=AVERAGE(0.2,0.27)"
